$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '286.55'
$ws.Range('E2').Value = '2.14%'
$ws.Range('G2').Value = '23'
$ws.Range('E3').Value = '3.98%'
$ws.Range('G3').Value = '23'
$ws.Range('D4').Value = '5.063'
$ws.Range('E4').Value = '4.54%'
$ws.Range('G4').Value = '23'
$ws.Range('E5').Value = '3.43%'
$ws.Range('G5').Value = '23'
$ws.Range('D6').Value = '7.392'
$ws.Range('E6').Value = '4.39%'
$ws.Range('G6').Value = '23'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '3.405'
$ws.Range('E7').Value = '3.03%'
$ws.Range('G7').Value = '23'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').Value = '1.366'
$ws.Range('E8').Value = '7.17%'
$ws.Range('G8').Value = '23'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9355'
$ws.Range('E9').Value = '3.82%'
$ws.Range('G9').Value = '23'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1578'
$ws.Range('E10').Value = '2.14%'
$ws.Range('G10').Value = '23'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '0.06715'
$ws.Range('E11').Value = '2.81%'
$ws.Range('G11').Value = '23'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.07602'
$ws.Range('E12').Value = '1.49%'
$ws.Range('G12').Value = '23'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.02935'
$ws.Range('E13').Value = '0.18%'
$ws.Range('G13').Value = '23'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.08993'
$ws.Range('E14').Value = '0.02%'
$ws.Range('G14').Value = '23'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001585'
$ws.Range('E15').Value = '-0.26%'
$ws.Range('G15').Value = '23'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').Value = '0.04504'
$ws.Range('E16').Value = '1.95%'
$ws.Range('G16').Value = '23'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').Value = '0.0006469'
$ws.Range('E17').Value = '0.62%'
$ws.Range('G17').Value = '23'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').Value = '0.006261'
$ws.Range('E18').Value = '3.42%'
$ws.Range('G18').Value = '23'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D19').Value = '3.446'
$ws.Range('E19').Value = '-1.20%'
$ws.Range('G19').Value = '23'
$ws.Range('D20').Value = '2.250'
$ws.Range('E20').Value = '1.15%'
$ws.Range('G20').Value = '23'
$ws.Range('E21').Value = '2.32%'
$ws.Range('G21').Value = '23'
$ws.Range('D22').Value = '0.1297'
$ws.Range('E22').Value = '-4.13%'
$ws.Range('G22').Value = '23'
$ws.Range('D23').Value = '4.097'
$ws.Range('E23').Value = '5.19%'
$ws.Range('G23').Value = '23'
$ws.Range('D24').Value = '0.1551'
$ws.Range('E24').Value = '3.30%'
$ws.Range('G24').Value = '23'
$ws.Range('D25').Value = '0.001179'
$ws.Range('G25').Value = '23'
$ws.Range('E26').Value = '-3.67%'
$ws.Range('G26').Value = '23'
$ws.Range('D27').Value = '0.0001248'
$ws.Range('E27').Value = '5.96%'
$ws.Range('G27').Value = '23'
$ws.Range('D28').Value = '0.0001616'
$ws.Range('E28').Value = '-2.39%'
$ws.Range('G28').Value = '23'
$ws.Range('G29').Value = '23'
$ws.Range('G30').Value = '23'
$ws.Range('G31').Value = '23'
$ws.Range('G32').Value = '23'
$ws.Range('G33').Value = '23'
$ws.Range('G34').Value = '23'
$ws.Range('G35').Value = '23'
$ws.Range('G36').Value = '23'
$ws.Range('G37').Value = '23'
$ws.Range('G38').Value = '23'
$ws.Range('G39').Value = '23'
$ws.Range('D40').Value = '0.04194'
$ws.Range('E40').Value = '2.59%'
$ws.Range('G40').Value = '23'
$ws.Range('D41').Value = '0.006716'
$ws.Range('E41').Value = '1.33%'
$ws.Range('G41').Value = '23'
$ws.Range('D42').Value = '0.1248'
$ws.Range('E42').Value = '-10.60%'
$ws.Range('G42').Value = '23'
$ws.Range('D43').Value = '0.002017'
$ws.Range('E43').Value = '-2.86%'
$ws.Range('G43').Value = '23'
$ws.Range('D44').Value = '0.01219'
$ws.Range('E44').Value = '10.51%'
$ws.Range('G44').Value = '23'
$ws.Range('D45').Value = '0.00005615'
$ws.Range('E45').Value = '1.18%'
$ws.Range('G45').Value = '23'
$ws.Range('E46').Value = '25.93%'
$ws.Range('G46').Value = '23'
$ws.Range('D47').Value = '0.01305'
$ws.Range('E47').Value = '-29.33%'
$ws.Range('G47').Value = '23'
$ws.Range('G48').Value = '23'
$ws.Range('G49').Value = '23'
$ws.Range('G50').Value = '23'
$ws.Range('G51').Value = '23'
